$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = 1
$ws.Range("B2").Value = "02A"
$ws.Range("C2").Value = 3

$ws.Range("A3").Value = 2
$ws.Range("B3").Value = "03A"
$ws.Range("C3").Value = 4

$ws.Range("A4").Value = 3
$ws.Range("B4").Value = "04A"
$ws.Range("C4").Value = 2

$ws.Range("A5").Value = 4
$ws.Range("B5").Value = "05A"
$ws.Range("C5").Value = 3

$ws.Range("A6").Value = 5
$ws.Range("B6").Value = "06A"
$ws.Range("C6").Value = 3

$ws.Range("A7").Value = 6
$ws.Range("B7").Value = "02B"
$ws.Range("C7").Value = 1

$ws.Range("A8").Value = 7
$ws.Range("B8").Value = "03B"
$ws.Range("C8").Value = 4

$ws.Range("A9").Value = 8
$ws.Range("B9").Value = "04B"
$ws.Range("C9").Value = 2

$ws.Range("A10").Value = 9
$ws.Range("B10").Value = "05B"
$ws.Range("C10").Value = 1

$ws.Range("A11").Value = 10
$ws.Range("B11").Value = "10B"
$ws.Range("C11").Value = 4

$ws.Range("A12").Value = 11
$ws.Range("B12").Value = "04C"
$ws.Range("C12").Value = 1

$ws.Range("A13").Value = 12
$ws.Range("B13").Value = "05C"
$ws.Range("C13").Value = 2

$ws.Range("A14").Value = 13
$ws.Range("B14").Value = "06C"
$ws.Range("C14").Value = 4

$ws.Range("A15").Value = 14
$ws.Range("B15").Value = "09C"
$ws.Range("C15").Value = 4

$ws.Range("A16").Value = 15
$ws.Range("B16").Value = "10C"
$ws.Range("C16").Value = 2

$ws.Range("A17").Value = 16
$ws.Range("B17").Value = "04D"
$ws.Range("C17").Value = 3

$ws.Range("A18").Value = 17
$ws.Range("B18").Value = "05D"
$ws.Range("C18").Value = 3

$ws.Range("A19").Value = 18
$ws.Range("B19").Value = "06D"
$ws.Range("C19").Value = 2

$ws.Range("A20").Value = 19
$ws.Range("B20").Value = "08D"
$ws.Range("C20").Value = 1

$ws.Range("A21").Value = 20
$ws.Range("B21").Value = "10D"
$ws.Range("C21").Value = 1

$ws.Range("A22").Value = 21
$ws.Range("B22").Value = "04E"
$ws.Range("C22").Value = 1

$ws.Range("A23").Value = 22
$ws.Range("B23").Value = "05E"
$ws.Range("C23").Value = 4

$ws.Range("A24").Value = 23
$ws.Range("B24").Value = "07E"
$ws.Range("C24").Value = 1

$ws.Range("A25").Value = 24
$ws.Range("B25").Value = "02F"
$ws.Range("C25").Value = 3

$ws.Range("A26").Value = 25
$ws.Range("B26").Value = "06F"
$ws.Range("C26").Value = 2

$ws.Range("A27").Value = 26
$ws.Range("B27").Value = "09F"
$ws.Range("C27").Value = 2

$ws.Range("A28").Value = 27
$ws.Range("B28").Value = "04G"
$ws.Range("C28").Value = 2

$ws.Range("A29").Value = 28
$ws.Range("B29").Value = "06G"
$ws.Range("C29").Value = 2

$ws.Range("A30").Value = 29
$ws.Range("B30").Value = "09G"
$ws.Range("C30").Value = 1

$ws.Range("A31").Value = 30
$ws.Range("B31").Value = "02H"
$ws.Range("C31").Value = 4

$ws.Range("A32").Value = 31
$ws.Range("B32").Value = "05H"
$ws.Range("C32").Value = 4

$ws.Range("A33").Value = 32
$ws.Range("B33").Value = "07H"
$ws.Range("C33").Value = 4

$ws.Range("A34").Value = 33
$ws.Range("B34").Value = "08H"
$ws.Range("C34").Value = 3
